$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row [15] (diff @1379)
$ws.Range("H15").Value = 182
$ws.Range("I15").Value = 182
$ws.Range("K15").Value = 546
$ws.Range("M15").Value = -377
# row [33] (diff @2282)
$ws.Range("H33").Value = 10417044
$ws.Range("I33").Value = 12346053
$ws.Range("J33").Value = 398
$ws.Range("K33").Value = 12346053
$ws.Range("L33").Value = 398
$ws.Range("M33").Value = -12345824
$ws.Range("N33").Value = -856
# row [40] (diff @2631)
$ws.Range("H40").Value = 2148.1667
$ws.Range("I40").Value = 3495.6
$ws.Range("J40").Value = 1629.9231
$ws.Range("K40").Value = 3495.6
$ws.Range("L40").Value = 1629.9231
$ws.Range("M40").Value = -3320.6
$ws.Range("N40").Value = -1979.9231
# row [101] (diff @5692)
$ws.Range("H101").Value = 57757
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 57757
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 173271
$ws.Range("N101").Value = -176515
$ws.Range("M101").ClearContents()
# row [103] (diff @5790)
$ws.Range("H103").Value = 398.5
$ws.Range("I103").Value = 399
$ws.Range("J103").Value = 398
$ws.Range("K103").Value = 1197
$ws.Range("L103").Value = 1194
$ws.Range("M103").Value = -611
$ws.Range("N103").Value = -2366
# row [113] (diff @6298)
$ws.Range("H113").Value = 2415.0588
$ws.Range("I113").Value = 1641.6666
$ws.Range("J113").Value = 2836.9092
$ws.Range("K113").Value = 1641.6666
$ws.Range("L113").Value = 2836.9092
$ws.Range("M113").Value = 1612.3334
$ws.Range("N113").Value = -9344.9092
# row [125] (diff @6892)
$ws.Range("H125").Value = 1607.2
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
# row [138] (diff @7547)
$ws.Range("H138").Value = 2885.58
$ws.Range("I138").Value = 1184.9678
$ws.Range("J138").Value = 3649.6233
$ws.Range("K138").Value = 3554.9034
$ws.Range("L138").Value = 10948.8699
$ws.Range("M138").Value = 1585.0966
$ws.Range("N138").Value = -21228.8699
# row [141] (diff @7691)
$ws.Range("H141").Value = 2950.1904
$ws.Range("I141").Value = 2576.5789
$ws.Range("J141").Value = 6499.5
$ws.Range("K141").Value = 7729.736699999999
$ws.Range("L141").Value = 19498.5
$ws.Range("M141").Value = -2549.736699999999
$ws.Range("N141").Value = -29858.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row [2] (diff @7837)
$ws.Range("H2").Value = 1295.5625
$ws.Range("I2").Value = 714.1429000000001
$ws.Range("J2").Value = 1747.7778
$ws.Range("K2").Value = 714.1429000000001
$ws.Range("L2").Value = 1747.7778
$ws.Range("M2").Value = -601.1429000000001
$ws.Range("N2").Value = -1973.7778
# row [32] (diff @9292)
$ws.Range("H32").Value = 55875.8
$ws.Range("I32").Value = 55214.176
$ws.Range("J32").Value = 63484.5
$ws.Range("K32").Value = 55214.176
$ws.Range("L32").Value = 63484.5
$ws.Range("M32").Value = -54927.176
$ws.Range("N32").Value = -64058.5
# row [45] (diff @9923)
$ws.Range("H45").Value = 893.35297
$ws.Range("I45").Value = 829.8461
$ws.Range("J45").Value = 1099.75
$ws.Range("K45").Value = 829.8461
$ws.Range("L45").Value = 1099.75
$ws.Range("M45").Value = -452.8461
$ws.Range("N45").Value = -1853.75
# row [63] (diff @10796)
$ws.Range("H63").Value = 1117089.5
$ws.Range("I63").Value = 3335335
$ws.Range("J63").Value = 7966.6665
$ws.Range("K63").Value = 3335335
$ws.Range("L63").Value = 7966.6665
$ws.Range("M63").Value = -3334649
$ws.Range("N63").Value = -9338.666499999999
# row [66] (diff @10946)
$ws.Range("H66").Value = 1117089.5
$ws.Range("I66").Value = 3335335
$ws.Range("J66").Value = 7966.6665
$ws.Range("K66").Value = 16676675
$ws.Range("L66").Value = 39833.3325
$ws.Range("M66").Value = -16673243
$ws.Range("N66").Value = -46697.3325
# row [116] (diff @13399)
$ws.Range("H116").Value = 1295.5625
$ws.Range("I116").Value = 714.1429000000001
$ws.Range("J116").Value = 1747.7778
$ws.Range("K116").Value = 714.1429000000001
$ws.Range("L116").Value = 1747.7778
$ws.Range("M116").Value = 1579.8571
$ws.Range("N116").Value = -6335.7778

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row [3] (diff @14822)
$ws.Range("H3").Value = 1295.5625
$ws.Range("I3").Value = 714.1429000000001
$ws.Range("J3").Value = 1747.7778
$ws.Range("K3").Value = 714.1429000000001
$ws.Range("L3").Value = 1747.7778
$ws.Range("M3").Value = -600.1429000000001
$ws.Range("N3").Value = -1975.7778
# row [20] (diff @15667)
$ws.Range("H20").Value = 2979.5518
$ws.Range("I20").Value = 3171.9583
$ws.Range("J20").Value = 2056
$ws.Range("K20").Value = 3171.9583
$ws.Range("L20").Value = 2056
$ws.Range("M20").Value = -2924.9583
$ws.Range("N20").Value = -2550
# row [107] (diff @19921)
$ws.Range("H107").Value = 2588.875
$ws.Range("I107").Value = 2285.1667
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 2285.1667
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = -365.1667000000002
$ws.Range("N107").Value = -7340

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row [94] (diff @26229)
$ws.Range("H94").Value = 125001430
$ws.Range("J94").Value = 1634.2858
$ws.Range("L94").Value = 1634.2858
$ws.Range("N94").Value = -2536.2858
# row [99] (diff @26474)
$ws.Range("H99").Value = 2072.7354
$ws.Range("I99").Value = 1722.2632
$ws.Range("J99").Value = 2516.6667
$ws.Range("K99").Value = 1722.2632
$ws.Range("L99").Value = 2516.6667
$ws.Range("M99").Value = -224.2632000000001
$ws.Range("N99").Value = -5512.6667
# row [126] (diff @27806)
$ws.Range("H126").Value = 2072.7354
$ws.Range("I126").Value = 1722.2632
$ws.Range("J126").Value = 2516.6667
$ws.Range("K126").Value = 5166.7896
$ws.Range("L126").Value = 7550.000100000001
$ws.Range("M126").Value = -2696.7896
$ws.Range("N126").Value = -12490.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row [131] (diff @35215)
$ws.Range("H131").Value = 734.22
$ws.Range("J131").Value = 805.1395
$ws.Range("L131").Value = 2415.4185
$ws.Range("N131").Value = -12495.4185

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row [43] (diff @37881)
$ws.Range("H43").Value = 1921.9
$ws.Range("I43").Value = 802.1111
$ws.Range("K43").Value = 802.1111
$ws.Range("M43").Value = -651.1111
# row [46] (diff @38025)
$ws.Range("H46").Value = 3633.3333
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# row [57] (diff @38567)
$ws.Range("H57").Value = 15233.333
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19640
# row [80] (diff @39685)
$ws.Range("H80").Value = 9750.5
$ws.Range("I80").Value = 7000.6665
$ws.Range("K80").Value = 7000.6665
$ws.Range("M80").Value = -6002.6665
# row [83] (diff @39832)
$ws.Range("H83").Value = 9750.5
$ws.Range("I83").Value = 7000.6665
$ws.Range("K83").Value = 35003.3325
$ws.Range("M83").Value = -30011.3325
# row [102] (diff @40751)
$ws.Range("H102").Value = 2343.2144
$ws.Range("I102").Value = 2687.889
$ws.Range("J102").Value = 1722.8
$ws.Range("K102").Value = 2687.889
$ws.Range("L102").Value = 1722.8
$ws.Range("M102").Value = -1065.889
$ws.Range("N102").Value = -4966.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row [46] (diff @44970)
$ws.Range("H46").Value = 1503.421
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 1536.9445
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1536.9445
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1912.9445
# row [55] (diff @45405)
$ws.Range("H55").Value = 284.05127
$ws.Range("I55").Value = 299.05554
$ws.Range("K55").Value = 299.05554
$ws.Range("M55").Value = -126.05554
# row [122] (diff @48670)
$ws.Range("H122").Value = 4438
$ws.Range("I122").Value = 6901.3335
$ws.Range("J122").Value = 2960
$ws.Range("K122").Value = 20704.0005
$ws.Range("L122").Value = 8880
$ws.Range("M122").Value = -18254.0005
$ws.Range("N122").Value = -13780

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row [5] (diff @49897)
$ws.Range("H5").Value = 37600000
$ws.Range("I5").Value = 37600000
$ws.Range("K5").Value = 37600000
$ws.Range("M5").Value = -37599888
# row [122] (diff @55516)
$ws.Range("H122").Value = 4435.9287
$ws.Range("I122").Value = 4810.4
$ws.Range("J122").Value = 3499.75
$ws.Range("K122").Value = 14431.2
$ws.Range("L122").Value = 10499.25
$ws.Range("M122").Value = -11981.2
$ws.Range("N122").Value = -15399.25

Write-Output "applied 196 value changes across 8 sheets"